$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings that Excel would
# otherwise coerce into numbers (dropping meaningful trailing zeros or flipping
# into scientific notation). Force them to Text format first so the literal
# string is preserved, then restore the original (default/"Normal") cell style
# so no visible formatting change is left behind.
$textForced = @('D20', 'D34', 'D40', 'D50')
foreach ($addr in $textForced) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.937.44'
$ws.Range('D3').Value = '1.895.36'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '0.7746'
$ws.Range('D6').Value = '243.94'
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('D9').Value = '25.83'
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').Value = '0.07373'
$ws.Range('E10').Value = '  +4.76%  '
$ws.Range('D11').Value = '0.08073'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '0.7738'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '5.511'
$ws.Range('E13').Value = '  +3.34%  '
$ws.Range('D14').Value = '1.897.34'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').Value = '94.35'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('E16').Value = '  +3.91%  '
$ws.Range('D17').Value = '29.927.17'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').Value = '247.63'
$ws.Range('E19').Value = '  +1.33%  '
$ws.Range('D20').Value = '0.000007852'
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('D21').Value = '8.179'
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('D22').Value = '2.154.98'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '0.1581'
$ws.Range('E25').Value = '  -4.58%  '
$ws.Range('D26').Value = '9.473'
$ws.Range('E26').Value = '  +1.45%  '
$ws.Range('D27').Value = '163.28'
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('D28').Value = '18.75'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').Value = '1.431'
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').Value = '4.476'
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('D33').Value = '0.05578'
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('D34').Value = '4.070'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = '1.243'
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('D36').Value = '0.7556'
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').Value = '1.006'
$ws.Range('E37').Value = '  +0.64%  '
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('D39').Value = '0.01934'
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('D40').Value = '2.790'
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '74.49'
$ws.Range('E41').Value = '  +2.55%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.4476'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = '1.108.41'
$ws.Range('E43').Value = '  +7.05%  '
$ws.Range('D44').Value = '6.006'
$ws.Range('E44').Value = '  +3.14%  '
$ws.Range('D45').Value = '0.8519'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = '1.895'
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('D48').Value = '102.58'
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.827'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '7.540'
$ws.Range('E50').Value = '  +1.33%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.094.15'
$ws.Range('E51').Value = '  +2.46%  '

foreach ($addr in $textForced) {
    $ws.Range($addr).Style = "Normal"
}

